# Replace every arithmetic "answer" cell in the single 20x5 table with its
# new expression/result, cell by cell (by row/column index) so that each
# Find.Execute only ever targets the exact cell intended - this matters
# because a few original expressions repeat verbatim in more than one
# cell (e.g. "35-7=28") but map to different replacements.
#
# Note: Find.Execute's "Replace" argument uses wdReplaceOne (1) rather than
# wdReplaceAll (2). Even though each call is scoped to a single cell's
# Range, wdReplaceAll was observed to replace every matching occurrence in
# the whole table, not just within that Range - which corrupted the other
# cell(s) sharing the same original text. wdReplaceOne replaces just the
# single occurrence found inside the cell's own Range, as intended.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$c = $t.Cell(1, 1).Range
$c.Find.Execute("74-26=48", $true, $false, $false, $false, $false, $true, 1, $false, "39+23=62", 1) | Out-Null
$c = $t.Cell(1, 2).Range
$c.Find.Execute("47+25=72", $true, $false, $false, $false, $false, $true, 1, $false, "16+77=93", 1) | Out-Null
$c = $t.Cell(1, 3).Range
$c.Find.Execute("54-46=8", $true, $false, $false, $false, $false, $true, 1, $false, "66-38=28", 1) | Out-Null
$c = $t.Cell(1, 4).Range
$c.Find.Execute("43-25=18", $true, $false, $false, $false, $false, $true, 1, $false, "6+75=81", 1) | Out-Null
$c = $t.Cell(1, 5).Range
$c.Find.Execute("19+13=32", $true, $false, $false, $false, $false, $true, 1, $false, "97-39=58", 1) | Out-Null
$c = $t.Cell(2, 1).Range
$c.Find.Execute("84-68=16", $true, $false, $false, $false, $false, $true, 1, $false, "70-14=56", 1) | Out-Null
$c = $t.Cell(2, 2).Range
$c.Find.Execute("17+35=52", $true, $false, $false, $false, $false, $true, 1, $false, "59+3=62", 1) | Out-Null
$c = $t.Cell(2, 3).Range
$c.Find.Execute("39+18=57", $true, $false, $false, $false, $false, $true, 1, $false, "18+78=96", 1) | Out-Null
$c = $t.Cell(2, 4).Range
$c.Find.Execute("84-27=57", $true, $false, $false, $false, $false, $true, 1, $false, "43-19=24", 1) | Out-Null
$c = $t.Cell(2, 5).Range
$c.Find.Execute("49+45=94", $true, $false, $false, $false, $false, $true, 1, $false, "54-7=47", 1) | Out-Null
$c = $t.Cell(3, 1).Range
$c.Find.Execute("73-15=58", $true, $false, $false, $false, $false, $true, 1, $false, "9+29=38", 1) | Out-Null
$c = $t.Cell(3, 2).Range
$c.Find.Execute("18+14=32", $true, $false, $false, $false, $false, $true, 1, $false, "35+57=92", 1) | Out-Null
$c = $t.Cell(3, 3).Range
$c.Find.Execute("91-29=62", $true, $false, $false, $false, $false, $true, 1, $false, "6+7=13", 1) | Out-Null
$c = $t.Cell(3, 4).Range
$c.Find.Execute("28-9=19", $true, $false, $false, $false, $false, $true, 1, $false, "18+46=64", 1) | Out-Null
$c = $t.Cell(3, 5).Range
$c.Find.Execute("71-65=6", $true, $false, $false, $false, $false, $true, 1, $false, "64+28=92", 1) | Out-Null
$c = $t.Cell(4, 1).Range
$c.Find.Execute("83-39=44", $true, $false, $false, $false, $false, $true, 1, $false, "28+27=55", 1) | Out-Null
$c = $t.Cell(4, 2).Range
$c.Find.Execute("44-29=15", $true, $false, $false, $false, $false, $true, 1, $false, "87+8=95", 1) | Out-Null
$c = $t.Cell(4, 3).Range
$c.Find.Execute("82-79=3", $true, $false, $false, $false, $false, $true, 1, $false, "7+47=54", 1) | Out-Null
$c = $t.Cell(4, 4).Range
$c.Find.Execute("16+15=31", $true, $false, $false, $false, $false, $true, 1, $false, "17+68=85", 1) | Out-Null
$c = $t.Cell(4, 5).Range
$c.Find.Execute("47+15=62", $true, $false, $false, $false, $false, $true, 1, $false, "11-5=6", 1) | Out-Null
$c = $t.Cell(5, 1).Range
$c.Find.Execute("26+19=45", $true, $false, $false, $false, $false, $true, 1, $false, "72-19=53", 1) | Out-Null
$c = $t.Cell(5, 2).Range
$c.Find.Execute("25+38=63", $true, $false, $false, $false, $false, $true, 1, $false, "28+44=72", 1) | Out-Null
$c = $t.Cell(5, 3).Range
$c.Find.Execute("57-18=39", $true, $false, $false, $false, $false, $true, 1, $false, "71-39=32", 1) | Out-Null
$c = $t.Cell(5, 4).Range
$c.Find.Execute("75-67=8", $true, $false, $false, $false, $false, $true, 1, $false, "43+28=71", 1) | Out-Null
$c = $t.Cell(5, 5).Range
$c.Find.Execute("26+15=41", $true, $false, $false, $false, $false, $true, 1, $false, "91-83=8", 1) | Out-Null
$c = $t.Cell(6, 1).Range
$c.Find.Execute("12+49=61", $true, $false, $false, $false, $false, $true, 1, $false, "74+19=93", 1) | Out-Null
$c = $t.Cell(6, 2).Range
$c.Find.Execute("61-55=6", $true, $false, $false, $false, $false, $true, 1, $false, "39+27=66", 1) | Out-Null
$c = $t.Cell(6, 3).Range
$c.Find.Execute("6+37=43", $true, $false, $false, $false, $false, $true, 1, $false, "71-63=8", 1) | Out-Null
$c = $t.Cell(6, 4).Range
$c.Find.Execute("64-19=45", $true, $false, $false, $false, $false, $true, 1, $false, "56+16=72", 1) | Out-Null
$c = $t.Cell(6, 5).Range
$c.Find.Execute("31-8=23", $true, $false, $false, $false, $false, $true, 1, $false, "27+57=84", 1) | Out-Null
$c = $t.Cell(7, 1).Range
$c.Find.Execute("17+44=61", $true, $false, $false, $false, $false, $true, 1, $false, "69+16=85", 1) | Out-Null
$c = $t.Cell(7, 2).Range
$c.Find.Execute("66+18=84", $true, $false, $false, $false, $false, $true, 1, $false, "62-28=34", 1) | Out-Null
$c = $t.Cell(7, 3).Range
$c.Find.Execute("43-9=34", $true, $false, $false, $false, $false, $true, 1, $false, "43-27=16", 1) | Out-Null
$c = $t.Cell(7, 4).Range
$c.Find.Execute("27-19=8", $true, $false, $false, $false, $false, $true, 1, $false, "23-19=4", 1) | Out-Null
$c = $t.Cell(7, 5).Range
$c.Find.Execute("33+59=92", $true, $false, $false, $false, $false, $true, 1, $false, "60-6=54", 1) | Out-Null
$c = $t.Cell(8, 1).Range
$c.Find.Execute("42+9=51", $true, $false, $false, $false, $false, $true, 1, $false, "16+38=54", 1) | Out-Null
$c = $t.Cell(8, 2).Range
$c.Find.Execute("16+57=73", $true, $false, $false, $false, $false, $true, 1, $false, "81-5=76", 1) | Out-Null
$c = $t.Cell(8, 3).Range
$c.Find.Execute("9+5=14", $true, $false, $false, $false, $false, $true, 1, $false, "90-77=13", 1) | Out-Null
$c = $t.Cell(8, 4).Range
$c.Find.Execute("19+65=84", $true, $false, $false, $false, $false, $true, 1, $false, "59+22=81", 1) | Out-Null
$c = $t.Cell(8, 5).Range
$c.Find.Execute("19+34=53", $true, $false, $false, $false, $false, $true, 1, $false, "44+28=72", 1) | Out-Null
$c = $t.Cell(9, 1).Range
$c.Find.Execute("69+9=78", $true, $false, $false, $false, $false, $true, 1, $false, "78+7=85", 1) | Out-Null
$c = $t.Cell(9, 2).Range
$c.Find.Execute("36+15=51", $true, $false, $false, $false, $false, $true, 1, $false, "53-7=46", 1) | Out-Null
$c = $t.Cell(9, 3).Range
$c.Find.Execute("22-9=13", $true, $false, $false, $false, $false, $true, 1, $false, "60-13=47", 1) | Out-Null
$c = $t.Cell(9, 4).Range
$c.Find.Execute("91-6=85", $true, $false, $false, $false, $false, $true, 1, $false, "40-19=21", 1) | Out-Null
$c = $t.Cell(9, 5).Range
$c.Find.Execute("28+24=52", $true, $false, $false, $false, $false, $true, 1, $false, "58+14=72", 1) | Out-Null
$c = $t.Cell(10, 1).Range
$c.Find.Execute("35-7=28", $true, $false, $false, $false, $false, $true, 1, $false, "91-42=49", 1) | Out-Null
$c = $t.Cell(10, 2).Range
$c.Find.Execute("34-8=26", $true, $false, $false, $false, $false, $true, 1, $false, "36+45=81", 1) | Out-Null
$c = $t.Cell(10, 3).Range
$c.Find.Execute("66+6=72", $true, $false, $false, $false, $false, $true, 1, $false, "7+74=81", 1) | Out-Null
$c = $t.Cell(10, 4).Range
$c.Find.Execute("29+63=92", $true, $false, $false, $false, $false, $true, 1, $false, "6+45=51", 1) | Out-Null
$c = $t.Cell(10, 5).Range
$c.Find.Execute("35-7=28", $true, $false, $false, $false, $false, $true, 1, $false, "62-58=4", 1) | Out-Null
$c = $t.Cell(11, 1).Range
$c.Find.Execute("40-3=37", $true, $false, $false, $false, $false, $true, 1, $false, "24+48=72", 1) | Out-Null
$c = $t.Cell(11, 2).Range
$c.Find.Execute("95-9=86", $true, $false, $false, $false, $false, $true, 1, $false, "7+29=36", 1) | Out-Null
$c = $t.Cell(11, 3).Range
$c.Find.Execute("38+46=84", $true, $false, $false, $false, $false, $true, 1, $false, "51-27=24", 1) | Out-Null
$c = $t.Cell(11, 4).Range
$c.Find.Execute("19+42=61", $true, $false, $false, $false, $false, $true, 1, $false, "16+25=41", 1) | Out-Null
$c = $t.Cell(11, 5).Range
$c.Find.Execute("67+8=75", $true, $false, $false, $false, $false, $true, 1, $false, "54-8=46", 1) | Out-Null
$c = $t.Cell(12, 1).Range
$c.Find.Execute("35-26=9", $true, $false, $false, $false, $false, $true, 1, $false, "38+53=91", 1) | Out-Null
$c = $t.Cell(12, 2).Range
$c.Find.Execute("91-24=67", $true, $false, $false, $false, $false, $true, 1, $false, "51-19=32", 1) | Out-Null
$c = $t.Cell(12, 3).Range
$c.Find.Execute("19+49=68", $true, $false, $false, $false, $false, $true, 1, $false, "80-15=65", 1) | Out-Null
$c = $t.Cell(12, 4).Range
$c.Find.Execute("65-39=26", $true, $false, $false, $false, $false, $true, 1, $false, "58+24=82", 1) | Out-Null
$c = $t.Cell(12, 5).Range
$c.Find.Execute("9+55=64", $true, $false, $false, $false, $false, $true, 1, $false, "37+8=45", 1) | Out-Null
$c = $t.Cell(13, 1).Range
$c.Find.Execute("54+39=93", $true, $false, $false, $false, $false, $true, 1, $false, "46-7=39", 1) | Out-Null
$c = $t.Cell(13, 2).Range
$c.Find.Execute("83-49=34", $true, $false, $false, $false, $false, $true, 1, $false, "82-9=73", 1) | Out-Null
$c = $t.Cell(13, 3).Range
$c.Find.Execute("41-38=3", $true, $false, $false, $false, $false, $true, 1, $false, "85+7=92", 1) | Out-Null
$c = $t.Cell(13, 4).Range
$c.Find.Execute("93-47=46", $true, $false, $false, $false, $false, $true, 1, $false, "50-48=2", 1) | Out-Null
$c = $t.Cell(13, 5).Range
$c.Find.Execute("71-43=28", $true, $false, $false, $false, $false, $true, 1, $false, "95-76=19", 1) | Out-Null
$c = $t.Cell(14, 1).Range
$c.Find.Execute("82-14=68", $true, $false, $false, $false, $false, $true, 1, $false, "95-89=6", 1) | Out-Null
$c = $t.Cell(14, 2).Range
$c.Find.Execute("92-18=74", $true, $false, $false, $false, $false, $true, 1, $false, "64-25=39", 1) | Out-Null
$c = $t.Cell(14, 3).Range
$c.Find.Execute("27+65=92", $true, $false, $false, $false, $false, $true, 1, $false, "29+3=32", 1) | Out-Null
$c = $t.Cell(14, 4).Range
$c.Find.Execute("56-17=39", $true, $false, $false, $false, $false, $true, 1, $false, "60-25=35", 1) | Out-Null
$c = $t.Cell(14, 5).Range
$c.Find.Execute("40-7=33", $true, $false, $false, $false, $false, $true, 1, $false, "39+4=43", 1) | Out-Null
$c = $t.Cell(15, 1).Range
$c.Find.Execute("84-58=26", $true, $false, $false, $false, $false, $true, 1, $false, "65+29=94", 1) | Out-Null
$c = $t.Cell(15, 2).Range
$c.Find.Execute("75-47=28", $true, $false, $false, $false, $false, $true, 1, $false, "6+89=95", 1) | Out-Null
$c = $t.Cell(15, 3).Range
$c.Find.Execute("72-66=6", $true, $false, $false, $false, $false, $true, 1, $false, "7+84=91", 1) | Out-Null
$c = $t.Cell(15, 4).Range
$c.Find.Execute("27+68=95", $true, $false, $false, $false, $false, $true, 1, $false, "87-69=18", 1) | Out-Null
$c = $t.Cell(15, 5).Range
$c.Find.Execute("80-23=57", $true, $false, $false, $false, $false, $true, 1, $false, "92-89=3", 1) | Out-Null
$c = $t.Cell(16, 1).Range
$c.Find.Execute("33-8=25", $true, $false, $false, $false, $false, $true, 1, $false, "6+68=74", 1) | Out-Null
$c = $t.Cell(16, 2).Range
$c.Find.Execute("58+3=61", $true, $false, $false, $false, $false, $true, 1, $false, "70-21=49", 1) | Out-Null
$c = $t.Cell(16, 3).Range
$c.Find.Execute("18+3=21", $true, $false, $false, $false, $false, $true, 1, $false, "67+26=93", 1) | Out-Null
$c = $t.Cell(16, 4).Range
$c.Find.Execute("91-37=54", $true, $false, $false, $false, $false, $true, 1, $false, "93-74=19", 1) | Out-Null
$c = $t.Cell(16, 5).Range
$c.Find.Execute("28+39=67", $true, $false, $false, $false, $false, $true, 1, $false, "40-14=26", 1) | Out-Null
$c = $t.Cell(17, 1).Range
$c.Find.Execute("39+46=85", $true, $false, $false, $false, $false, $true, 1, $false, "38+5=43", 1) | Out-Null
$c = $t.Cell(17, 2).Range
$c.Find.Execute("42-5=37", $true, $false, $false, $false, $false, $true, 1, $false, "32+49=81", 1) | Out-Null
$c = $t.Cell(17, 3).Range
$c.Find.Execute("5+59=64", $true, $false, $false, $false, $false, $true, 1, $false, "94-7=87", 1) | Out-Null
$c = $t.Cell(17, 4).Range
$c.Find.Execute("68+25=93", $true, $false, $false, $false, $false, $true, 1, $false, "71-18=53", 1) | Out-Null
$c = $t.Cell(17, 5).Range
$c.Find.Execute("60-5=55", $true, $false, $false, $false, $false, $true, 1, $false, "9+84=93", 1) | Out-Null
$c = $t.Cell(18, 1).Range
$c.Find.Execute("58+25=83", $true, $false, $false, $false, $false, $true, 1, $false, "17+18=35", 1) | Out-Null
$c = $t.Cell(18, 2).Range
$c.Find.Execute("17+69=86", $true, $false, $false, $false, $false, $true, 1, $false, "9+82=91", 1) | Out-Null
$c = $t.Cell(18, 3).Range
$c.Find.Execute("50-38=12", $true, $false, $false, $false, $false, $true, 1, $false, "12-3=9", 1) | Out-Null
$c = $t.Cell(18, 4).Range
$c.Find.Execute("48+19=67", $true, $false, $false, $false, $false, $true, 1, $false, "19+3=22", 1) | Out-Null
$c = $t.Cell(18, 5).Range
$c.Find.Execute("72-37=35", $true, $false, $false, $false, $false, $true, 1, $false, "17+9=26", 1) | Out-Null
$c = $t.Cell(19, 1).Range
$c.Find.Execute("52-38=14", $true, $false, $false, $false, $false, $true, 1, $false, "50-39=11", 1) | Out-Null
$c = $t.Cell(19, 2).Range
$c.Find.Execute("33+9=42", $true, $false, $false, $false, $false, $true, 1, $false, "77+17=94", 1) | Out-Null
$c = $t.Cell(19, 3).Range
$c.Find.Execute("82-64=18", $true, $false, $false, $false, $false, $true, 1, $false, "44+49=93", 1) | Out-Null
$c = $t.Cell(19, 4).Range
$c.Find.Execute("64+17=81", $true, $false, $false, $false, $false, $true, 1, $false, "30-3=27", 1) | Out-Null
$c = $t.Cell(19, 5).Range
$c.Find.Execute("78-29=49", $true, $false, $false, $false, $false, $true, 1, $false, "12-5=7", 1) | Out-Null
$c = $t.Cell(20, 1).Range
$c.Find.Execute("15-6=9", $true, $false, $false, $false, $false, $true, 1, $false, "9+22=31", 1) | Out-Null
$c = $t.Cell(20, 2).Range
$c.Find.Execute("33+8=41", $true, $false, $false, $false, $false, $true, 1, $false, "6+28=34", 1) | Out-Null
$c = $t.Cell(20, 3).Range
$c.Find.Execute("60-57=3", $true, $false, $false, $false, $false, $true, 1, $false, "94-79=15", 1) | Out-Null
$c = $t.Cell(20, 4).Range
$c.Find.Execute("60-1=59", $true, $false, $false, $false, $false, $true, 1, $false, "33-18=15", 1) | Out-Null
$c = $t.Cell(20, 5).Range
$c.Find.Execute("92-55=37", $true, $false, $false, $false, $false, $true, 1, $false, "96-77=19", 1) | Out-Null

Write-Output "done"
